# Update "想去人数" (want-to-go count) and "最低票价" (min price) figures
# across the 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types)
# worksheets, per the latest scrape refresh.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1891
$ws.Range("F3").Value = 1519
$ws.Range("F4").Value = 886
$ws.Range("F5").Value = 780
$ws.Range("F6").Value = 13361
$ws.Range("F7").Value = 13229
$ws.Range("G8").Value = 60
$ws.Range("F9").Value = 778
$ws.Range("G9").Value = 55
$ws.Range("F11").Value = 565
$ws.Range("F13").Value = 682
$ws.Range("G13").Value = 60
$ws.Range("F14").Value = 2099
$ws.Range("F19").Value = 398
$ws.Range("F20").Value = 266
$ws.Range("F22").Value = 425

# --- 演出 (sheet 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 132
$ws.Range("F6").Value = 62

# --- 全部类型 (sheet 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1891
$ws.Range("F4").Value = 1519
$ws.Range("F5").Value = 886
$ws.Range("F7").Value = 780
$ws.Range("F8").Value = 13361
$ws.Range("F9").Value = 13229
$ws.Range("G10").Value = 60
$ws.Range("F11").Value = 778
$ws.Range("G11").Value = 55
$ws.Range("F13").Value = 565
$ws.Range("F15").Value = 682
$ws.Range("G15").Value = 60
$ws.Range("F18").Value = 2099
$ws.Range("F22").Value = 132
$ws.Range("F24").Value = 62
$ws.Range("F26").Value = 398
$ws.Range("F27").Value = 266
$ws.Range("F29").Value = 425
